# Adds a new "Gastos por m2" column (O, merged with P for the header)
# to the 2021 forest-fire sheet, with a style carrying a thousands-style
# number format + left alignment for the data cells, and a centered
# style for the merged header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style the whole data range (O2:O251) first: number format "#,##0" +
# left alignment. This becomes the shared cellXf used by every data cell
# in the new column, including the still-blank O2.
$dataRange = $ws.Range("O2:O251")
$dataRange.NumberFormat = "#,##0"
$dataRange.HorizontalAlignment = -4131

# Fill the two blocks of placeholder values (space vs double-space) —
# header cell is populated in between so the new shared strings land in
# the same order as the source edit: " ", "Gastos por m2", "  ".
$ws.Range("O3:O156").Value2 = " "
$ws.Range("O1").Value2 = "Gastos por m2"
$ws.Range("O157:O251").Value2 = "  "

# Header: centered alignment, merged across O1:P1.
$headerRange = $ws.Range("O1:P1")
$headerRange.HorizontalAlignment = -4108
$headerRange.Merge()

# Match the saved selection/scroll position from the source file.
$ws.Range("O3").Select()
